$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price record (Superior Seedless, Región de O'Higgins, 2023-01-25)
# was inserted as row 473, pushing every subsequent row down by one
# (old row 473 -> new row 474, ..., old row 574 -> new row 575).
$ws.Rows.Item(473).Insert()

$ws.Range("A473").Value = 5
$ws.Range("B473").Value = "Macroferia Regional de Talca"
$ws.Range("C473").Value = "Maule"
$ws.Range("D473").Value = 44951
$ws.Range("E473").Value = 7
$ws.Range("F473").Value = "Fruta"
$ws.Range("G473").Value = 100109
$ws.Range("H473").Value = "Uva"
$ws.Range("I473").Value = 100109001
$ws.Range("J473").Value = "Uva"
$ws.Range("K473").Value = "Superior Seedless"
$ws.Range("L473").Value = "Primera"
$ws.Range("M473").Value = 240
$ws.Range("N473").Value = 12000
$ws.Range("O473").Value = 12000
$ws.Range("P473").Value = 12000
$ws.Range("Q473").Value = "$/bandeja 18 kilos"
$ws.Range("R473").Value = "Región de O'Higgins"
$ws.Range("S473").Value = 667
$ws.Range("T473").Value = 18
